$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1449.8334
$ws.Range("J32").Value = 1499.8
$ws.Range("L32").Value = 1499.8
$ws.Range("N32").Value = -2151.8
$ws.Range("H116").Value = 2160
$ws.Range("H132").Value = 3862987
$ws.Range("I132").Value = 4763889
$ws.Range("K132").Value = 14291667
$ws.Range("M132").Value = -14289137
$ws.Range("H137").Value = 2252.818
$ws.Range("I137").Value = 2078.1
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 6234.299999999999
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -3684.299999999999
$ws.Range("N137").Value = -17100
$ws.Range("H141").Value = 2638.48
$ws.Range("I141").Value = 1559.6154
$ws.Range("J141").Value = 3807.25
$ws.Range("K141").Value = 4678.8462
$ws.Range("L141").Value = 11421.75
$ws.Range("M141").Value = 501.1538
$ws.Range("N141").Value = -21781.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1390.8235
$ws.Range("I2").Value = 754.7
$ws.Range("K2").Value = 754.7
$ws.Range("M2").Value = -641.7
$ws.Range("H45").Value = 924.04
$ws.Range("I45").Value = 1244
$ws.Range("J45").Value = 799.6111
$ws.Range("K45").Value = 1244
$ws.Range("L45").Value = 799.6111
$ws.Range("M45").Value = -867
$ws.Range("N45").Value = -1553.6111
$ws.Range("I64").Value = 41111
$ws.Range("K64").Value = 41111
$ws.Range("M64").Value = -40863
$ws.Range("I67").Value = 41111
$ws.Range("K67").Value = 41111
$ws.Range("M67").Value = -40253
$ws.Range("H116").Value = 1390.8235
$ws.Range("I116").Value = 754.7
$ws.Range("K116").Value = 754.7
$ws.Range("M116").Value = 1539.3
$ws.Range("H122").Value = 2174.8572
$ws.Range("I122").Value = 1258
$ws.Range("J122").Value = 2862.5
$ws.Range("K122").Value = 3774
$ws.Range("L122").Value = 8587.5
$ws.Range("M122").Value = -1324
$ws.Range("N122").Value = -13487.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1390.8235
$ws.Range("I3").Value = 754.7
$ws.Range("K3").Value = 754.7
$ws.Range("M3").Value = -640.7
$ws.Range("H20").Value = 3115.7727
$ws.Range("I20").Value = 3556.5833
$ws.Range("J20").Value = 2586.8
$ws.Range("K20").Value = 3556.5833
$ws.Range("L20").Value = 2586.8
$ws.Range("M20").Value = -3309.5833
$ws.Range("N20").Value = -3080.8
$ws.Range("H134").Value = 22322.041
$ws.Range("I134").Value = 31055.883
$ws.Range("K134").Value = 93167.649
$ws.Range("M134").Value = -90632.649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8337103
$ws.Range("I31").Value = 4146.143
$ws.Range("K31").Value = 4146.143
$ws.Range("M31").Value = -3851.143
$ws.Range("H34").Value = 8337103
$ws.Range("I34").Value = 4146.143
$ws.Range("K34").Value = 4146.143
$ws.Range("M34").Value = -3944.143
$ws.Range("H99").Value = 1772.7273
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 2307.6924
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 2307.6924
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -5303.6924
$ws.Range("H107").Value = 583.7742
$ws.Range("I107").Value = 535.6
$ws.Range("K107").Value = 535.6
$ws.Range("M107").Value = 1384.4
$ws.Range("H126").Value = 1772.7273
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 2307.6924
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 6923.0772
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -11863.0772
$ws.Range("H132").Value = 2644.28
$ws.Range("I132").Value = 1810.5
$ws.Range("K132").Value = 5431.5
$ws.Range("M132").Value = -2901.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 23205.8
$ws.Range("J101").Value = 23205.8
$ws.Range("L101").Value = 69617.39999999999
$ws.Range("N101").Value = -74485.39999999999
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 15981.818
$ws.Range("I110").Value = 2450
$ws.Range("J110").Value = 23714.285
$ws.Range("K110").Value = 7350
$ws.Range("L110").Value = 71142.855
$ws.Range("M110").Value = -3260
$ws.Range("N110").Value = -79322.855
$ws.Range("H113").Value = 514.29785
$ws.Range("I113").Value = 506.6154
$ws.Range("J113").Value = 523.8095
$ws.Range("K113").Value = 1519.8462
$ws.Range("L113").Value = 1571.4285
$ws.Range("M113").Value = 650.1538
$ws.Range("N113").Value = -5911.4285
$ws.Range("H114").Value = 3212.2144
$ws.Range("I114").Value = 2041.8
$ws.Range("J114").Value = 3862.4443
$ws.Range("K114").Value = 6125.4
$ws.Range("L114").Value = 11587.3329
$ws.Range("M114").Value = -2871.4
$ws.Range("N114").Value = -18095.3329
$ws.Range("H122").Value = 1574.9445
$ws.Range("I122").Value = 1495.1818
$ws.Range("K122").Value = 13456.6362
$ws.Range("M122").Value = -11006.6362
$ws.Range("H131").Value = 1062675.8
$ws.Range("I131").Value = 14894.286
$ws.Range("J131").Value = 1429399.2
$ws.Range("K131").Value = 44682.858
$ws.Range("L131").Value = 4288197.6
$ws.Range("M131").Value = -39642.858
$ws.Range("N131").Value = -4298277.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 38088.1
$ws.Range("J42").Value = 38088.1
$ws.Range("L42").Value = 38088.1
$ws.Range("N42").Value = -39058.1
$ws.Range("H111").Value = 5261
$ws.Range("J111").Value = 5261
$ws.Range("L111").Value = 5261
$ws.Range("N111").Value = -11395
$ws.Range("H115").Value = 38088.1
$ws.Range("J115").Value = 38088.1
$ws.Range("L115").Value = 38088.1
$ws.Range("N115").Value = -40438.1
$ws.Range("H118").Value = 13000
$ws.Range("J118").Value = 13000
$ws.Range("L118").Value = 13000
$ws.Range("N118").Value = -16314
$ws.Range("H122").Value = 1741
$ws.Range("I122").Value = 1575.8077
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 4727.4231
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -2277.4231
$ws.Range("N122").Value = -12700
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2137.3333
$ws.Range("I7").Value = 1840
$ws.Range("J7").Value = 2397.5
$ws.Range("K7").Value = 1840
$ws.Range("L7").Value = 2397.5
$ws.Range("M7").Value = -1728
$ws.Range("N7").Value = -2621.5
$ws.Range("H40").Value = 3500
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H126").Value = 2137.3333
$ws.Range("I126").Value = 1840
$ws.Range("J126").Value = 2397.5
$ws.Range("K126").Value = 5520
$ws.Range("L126").Value = 7192.5
$ws.Range("M126").Value = -3050
$ws.Range("N126").Value = -12132.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 52499.75
$ws.Range("J16").Value = 52499.75
$ws.Range("L16").Value = 52499.75
$ws.Range("N16").Value = -53083.75
$ws.Range("H81").Value = 3480
$ws.Range("I81").Value = 3466.6667
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 6933.3334
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = -5872.3334
$ws.Range("N81").Value = -9122
$ws.Range("H84").Value = 3480
$ws.Range("I84").Value = 3466.6667
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 34666.667
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = -29362.667
$ws.Range("N84").Value = -45608
$ws.Range("H126").Value = 6282.278
$ws.Range("I126").Value = 7805.4287
$ws.Range("J126").Value = 951.25
$ws.Range("K126").Value = 23416.2861
$ws.Range("L126").Value = 2853.75
$ws.Range("M126").Value = -20946.2861
$ws.Range("N126").Value = -7793.75
